$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2.242910268107096
$ws.Range("C2").Value = 0.6311885399363177
$ws.Range("E2").Value = 0.4158060904646845
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.00242407716902464
$ws.Range("N2").Value = 1.47328817965743
$ws.Range("B3").Value = 2.007184378090869
$ws.Range("C3").Value = 0.5525089403576544
$ws.Range("E3").Value = 0.362019698825975
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002432169680867075
$ws.Range("N3").Value = 1.462113483797097
$ws.Range("B4").Value = 1.863824805689944
$ws.Range("C4").Value = 0.5044814892216323
$ws.Range("E4").Value = 0.3291912763668421
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.00243737669552322
$ws.Range("N4").Value = 1.455779937646582
$ws.Range("B5").Value = 1.80573637306469
$ws.Range("C5").Value = 0.4849755974774439
$ws.Range("E5").Value = 0.3158581801401539
$ws.Range("F5").Value = 0.3390132514326325
$ws.Range("G5").Value = 0.002439558792655624
$ws.Range("N5").Value = 1.453328841223126
$ws.Range("B6").Value = 1.796110452037283
$ws.Range("C6").Value = 0.4817404744861733
$ws.Range("E6").Value = 0.3136467972641412
$ws.Range("F6").Value = 0.3366681778241372
$ws.Range("G6").Value = 0.002439924772921491
$ws.Range("N6").Value = 1.452929600530609
$ws.Range("B7").Value = 1.863040080462895
$ws.Range("C7").Value = 0.504218166867247
$ws.Range("E7").Value = 0.3290112865013697
$ws.Range("F7").Value = 0.3529483938368969
$ws.Range("G7").Value = 0.002437405880020995
$ws.Range("N7").Value = 1.455746359158738
$ws.Range("B8").Value = 2.161337265704276
$ws.Range("C8").Value = 0.6039980335335713
$ws.Range("E8").Value = 0.3972168606075144
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.00242681821776923
$ws.Range("N8").Value = 1.469323884945325
$ws.Range("B9").Value = 2.757920240944884
$ws.Range("C9").Value = 0.8021555606492257
$ws.Range("E9").Value = 0.5327497143518514
$ws.Range("F9").Value = 0.5661985755042025
$ws.Range("G9").Value = 0.002407931669278663
$ws.Range("N9").Value = 1.500268038415214
$ws.Range("B10").Value = 3.204341786139707
$ws.Range("C10").Value = 0.9496269551769956
$ws.Range("E10").Value = 0.6337339265794242
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002395179452167902
$ws.Range("N10").Value = 1.52582192964627
$ws.Range("B11").Value = 3.409423874693289
$ws.Range("C11").Value = 1.017207688569499
$ws.Range("E11").Value = 0.6800526058525094
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002389617830170714
$ws.Range("N11").Value = 1.538100673397508
$ws.Range("B12").Value = 3.48738982487464
$ws.Range("C12").Value = 1.042876692757204
$ws.Range("E12").Value = 0.6976529604841772
$ws.Range("F12").Value = 0.7356546913087669
$ws.Range("G12").Value = 0.002387545876413383
$ws.Range("N12").Value = 1.542847813480705
$ws.Range("B13").Value = 3.470584571606537
$ws.Range("C13").Value = 1.037344855597951
$ws.Range("E13").Value = 0.6938596257631389
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002387990596671201
$ws.Range("N13").Value = 1.541821043315906
$ws.Range("B14").Value = 3.415831947287302
$ws.Range("C14").Value = 1.019317899299381
$ws.Range("E14").Value = 0.6814993519042361
$ws.Range("F14").Value = 0.7191683204515442
$ws.Range("G14").Value = 0.002389446687452271
$ws.Range("N14").Value = 1.538489250867599
$ws.Range("B15").Value = 3.382334748440485
$ws.Range("C15").Value = 1.008286174138732
$ws.Range("E15").Value = 0.6739363772747566
$ws.Range("F15").Value = 0.7114413442032657
$ws.Range("G15").Value = 0.002390343017841352
$ws.Range("N15").Value = 1.536461223817298
$ws.Range("B16").Value = 3.190981000404236
$ws.Range("C16").Value = 0.9452209064223211
$ws.Range("E16").Value = 0.6307150366178007
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002395547707982758
$ws.Range("N16").Value = 1.525032935238301
$ws.Range("B17").Value = 3.074117127045156
$ws.Range("C17").Value = 0.9066637692233712
$ws.Range("E17").Value = 0.6043016352017929
$ws.Range("F17").Value = 0.6400460337216174
$ws.Range("G17").Value = 0.00239880171493173
$ws.Range("N17").Value = 1.518191974814272
$ws.Range("B18").Value = 3.007087521757057
$ws.Range("C18").Value = 0.8845329435203553
$ws.Range("E18").Value = 0.5891447999221668
$ws.Range("F18").Value = 0.6244449056557215
$ws.Range("G18").Value = 0.002400695888542229
$ws.Range("N18").Value = 1.514318600326305
$ws.Range("B19").Value = 2.984424101454351
$ws.Range("C19").Value = 0.8770475717351474
$ws.Range("E19").Value = 0.5840188725112654
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002401341106540452
$ws.Range("N19").Value = 1.513017579999911
$ws.Range("B20").Value = 3.086537962853072
$ws.Range("C20").Value = 0.9107634127544202
$ws.Range("E20").Value = 0.6071096762427715
$ws.Range("F20").Value = 0.6429339538360921
$ws.Range("G20").Value = 0.002398452988234602
$ws.Range("N20").Value = 1.518913831448344
$ws.Range("B21").Value = 3.431905675771077
$ws.Range("C21").Value = 1.024610693828947
$ws.Range("E21").Value = 0.6851281758088277
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002389018075002979
$ws.Range("N21").Value = 1.539465207076148
$ws.Range("B22").Value = 3.659414627157389
$ws.Range("C22").Value = 1.099472183521755
$ws.Range("E22").Value = 0.7364730819928269
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002383050494854145
$ws.Range("N22").Value = 1.553466335860037
$ws.Range("B23").Value = 3.537818890525841
$ws.Range("C23").Value = 1.059473324968621
$ws.Range("E23").Value = 0.7090348760984426
$ws.Range("F23").Value = 0.7472568307916134
$ws.Range("G23").Value = 0.002386217430972703
$ws.Range("N23").Value = 1.5459404125898
$ws.Range("B24").Value = 3.080922015330543
$ws.Range("C24").Value = 0.9089098516876106
$ws.Range("E24").Value = 0.6058400722986477
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002398610574551041
$ws.Range("N24").Value = 1.518587294725734
$ws.Range("B25").Value = 2.595180027534525
$ws.Range("C25").Value = 0.7482448054573752
$ws.Range("E25").Value = 0.4958609936548584
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.00241284220073668
$ws.Range("N25").Value = 1.49141544874891
